$d = $word.ActiveDocument

# Paragraph 1: "Студенту ${group} ${student}"
#   -> "Студенту " + MERGEFIELD student.group + " " + MERGEFIELD student.name
$para1Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2CFF1BC6" w14:textId="32BA4338" w:rsidR="0069296A" w:rsidRPr="00037CA0" w:rsidRDefault="00E1456F" w:rsidP="00E1456F"><w:pPr><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">Студенту </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:instrText xml:space="preserve"> MERGEFIELD  ${student.group}  \* MERGEFORMAT </w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:t>«${student.group}»</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:instrText>MERGEFIELD</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:instrText xml:space="preserve">  ${</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:instrText>student</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:instrText>.</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:instrText>name</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:instrText xml:space="preserve">}  \* </w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:instrText>MERGEFORMAT</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:t>«${</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>student</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>name</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:t>}»</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p>
'@

$rng1 = $d.Content
$found1 = $rng1.Find.Execute('${group} ${student}', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $found1) {
    throw 'Could not find placeholder ${group} ${student}'
}
$paraRng1 = $rng1.Paragraphs(1).Range
$paraRng1.InsertXML($para1Xml)

# Paragraph 2: Тема проекта " ${projectTheme} "
#   -> Тема проекта " " + MERGEFIELD topic.title + " ""
$para2Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2F9549CC" w14:textId="6C31EB10" w:rsidR="00B224F5" w:rsidRPr="00562918" w:rsidRDefault="00B224F5" w:rsidP="00B224F5"><w:pPr><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r w:rsidRPr="00B224F5"><w:rPr><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">1. </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">Тема проекта </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">" </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:instrText xml:space="preserve"> MERGEFIELD  ${topic.title}  \* MERGEFORMAT </w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:t>«${topic.title}»</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> "</w:t></w:r></w:p>
'@

$rng2 = $d.Content
$found2 = $rng2.Find.Execute('" ${projectTheme} "', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $found2) {
    throw 'Could not find placeholder " ${projectTheme} "'
}
$paraRng2 = $rng2.Paragraphs(1).Range
$paraRng2.InsertXML($para2Xml)

Write-Output 'done'
